$wb = $excel.ActiveWorkbook

# --- Update header text (shared strings) on both sheets that carry it ---
$ws1 = $wb.Worksheets.Item("Monthly Data")
$ws2 = $wb.Worksheets.Item("Annual Data")

$ws1.Range("A2").Value = "August 2020 Monthly Energy Review"
$ws1.Range("A6").Value = "Release Date: August 26, 2020"
$ws1.Range("A7").Value = "Next Update: September 24, 2020"

$ws2.Range("A2").Value = "August 2020 Monthly Energy Review"
$ws2.Range("A6").Value = "Release Date: August 26, 2020"
$ws2.Range("A7").Value = "Next Update: September 24, 2020"

# --- Revise monthly figures for Jan 2019 - Apr 2020 (rows 565-580) ---
# Row 565: 2=1204.702, 3=2545.862, 4=699.184, 5=1822.769, 6=2060.107, 7=2850.517
$ws1.Cells.Item(565,2).Value = 1204.702
$ws1.Cells.Item(565,3).Value = 2545.862
$ws1.Cells.Item(565,4).Value = 699.184
$ws1.Cells.Item(565,5).Value = 1822.769
$ws1.Cells.Item(565,6).Value = 2060.107
$ws1.Cells.Item(565,7).Value = 2850.517

# Row 566: 2=1022.817, 3=2158.604, 4=597.855, 5=1589.997, 6=1795.555, 7=2504.562
$ws1.Cells.Item(566,2).Value = 1022.817
$ws1.Cells.Item(566,3).Value = 2158.604
$ws1.Cells.Item(566,4).Value = 597.855
$ws1.Cells.Item(566,5).Value = 1589.997
$ws1.Cells.Item(566,6).Value = 1795.555
$ws1.Cells.Item(566,7).Value = 2504.562

# Row 567: 2=897.415, 3=2005.165, 4=551.019, 5=1606.905, 6=1946.862, 7=2709.452
$ws1.Cells.Item(567,2).Value = 897.415
$ws1.Cells.Item(567,3).Value = 2005.165
$ws1.Cells.Item(567,4).Value = 551.019
$ws1.Cells.Item(567,5).Value = 1606.905
$ws1.Cells.Item(567,6).Value = 1946.862
$ws1.Cells.Item(567,7).Value = 2709.452

# Row 568: 2=480.613, 3=1367.721, 4=346.036, 6=1839.092, 7=2593.416
$ws1.Cells.Item(568,2).Value = 480.613
$ws1.Cells.Item(568,3).Value = 1367.721
$ws1.Cells.Item(568,4).Value = 346.036
$ws1.Cells.Item(568,6).Value = 1839.092
$ws1.Cells.Item(568,7).Value = 2593.416

# Row 569: 2=347.638, 3=1362.475, 4=276.301, 5=1403.788, 6=1909.009, 7=2729.102
$ws1.Cells.Item(569,2).Value = 347.638
$ws1.Cells.Item(569,3).Value = 1362.475
$ws1.Cells.Item(569,4).Value = 276.301
$ws1.Cells.Item(569,5).Value = 1403.788
$ws1.Cells.Item(569,6).Value = 1909.009
$ws1.Cells.Item(569,7).Value = 2729.102

# Row 570: 2=246.372, 3=1448.276, 4=228.094, 5=1387.955, 6=1839.603, 7=2650.317
$ws1.Cells.Item(570,2).Value = 246.372
$ws1.Cells.Item(570,3).Value = 1448.276
$ws1.Cells.Item(570,4).Value = 228.094
$ws1.Cells.Item(570,5).Value = 1387.955
$ws1.Cells.Item(570,6).Value = 1839.603
$ws1.Cells.Item(570,7).Value = 2650.317

# Row 571: 6=1919.991
$ws1.Cells.Item(571,6).Value = 1919.991

# Row 572: 2=229.622, 3=1698.773, 4=234.265, 5=1512.359, 6=1949.946, 7=2798.191
$ws1.Cells.Item(572,2).Value = 229.622
$ws1.Cells.Item(572,3).Value = 1698.773
$ws1.Cells.Item(572,4).Value = 234.265
$ws1.Cells.Item(572,5).Value = 1512.359
$ws1.Cells.Item(572,6).Value = 1949.946
$ws1.Cells.Item(572,7).Value = 2798.191

# Row 573: 4=222.988, 6=1891.519
$ws1.Cells.Item(573,4).Value = 222.988
$ws1.Cells.Item(573,6).Value = 1891.519

# Row 574: 2=373.684, 3=1395.005, 4=310.771, 5=1396.007, 6=1974.967, 7=2733.478, 12=7958.937
$ws1.Cells.Item(574,2).Value = 373.684
$ws1.Cells.Item(574,3).Value = 1395.005
$ws1.Cells.Item(574,4).Value = 310.771
$ws1.Cells.Item(574,5).Value = 1396.007
$ws1.Cells.Item(574,6).Value = 1974.967
$ws1.Cells.Item(574,7).Value = 2733.478
$ws1.Cells.Item(574,12).Value = 7958.937

# Row 575: 2=782.414, 3=1811.48, 4=499.314, 5=1534.858, 6=1954.748, 7=2718.933
$ws1.Cells.Item(575,2).Value = 782.414
$ws1.Cells.Item(575,3).Value = 1811.48
$ws1.Cells.Item(575,4).Value = 499.314
$ws1.Cells.Item(575,5).Value = 1534.858
$ws1.Cells.Item(575,6).Value = 1954.748
$ws1.Cells.Item(575,7).Value = 2718.933

# Row 576: 2=980.091, 3=2172.512, 4=588.342, 5=1647.856, 6=2002.966, 7=2755.53
$ws1.Cells.Item(576,2).Value = 980.091
$ws1.Cells.Item(576,3).Value = 2172.512
$ws1.Cells.Item(576,4).Value = 588.342
$ws1.Cells.Item(576,5).Value = 1647.856
$ws1.Cells.Item(576,6).Value = 2002.966
$ws1.Cells.Item(576,7).Value = 2755.53

# Row 577: 2=1038.384, 3=2247.38, 4=618.663, 5=1671.164, 6=2004.026, 7=2760.16, 9=2277.144, 10=3024.613, 12=8951.407
$ws1.Cells.Item(577,2).Value = 1038.384
$ws1.Cells.Item(577,3).Value = 2247.38
$ws1.Cells.Item(577,4).Value = 618.663
$ws1.Cells.Item(577,5).Value = 1671.164
$ws1.Cells.Item(577,6).Value = 2004.026
$ws1.Cells.Item(577,7).Value = 2760.16
$ws1.Cells.Item(577,9).Value = 2277.144
$ws1.Cells.Item(577,10).Value = 3024.613
$ws1.Cells.Item(577,12).Value = 8951.407

# Row 578: 2=931.934, 3=2018.306, 4=563.265, 5=1553.334, 6=1867.119, 7=2600.919, 9=2139.218, 10=2816.266, 12=8305.8
$ws1.Cells.Item(578,2).Value = 931.934
$ws1.Cells.Item(578,3).Value = 2018.306
$ws1.Cells.Item(578,4).Value = 563.265
$ws1.Cells.Item(578,5).Value = 1553.334
$ws1.Cells.Item(578,6).Value = 1867.119
$ws1.Cells.Item(578,7).Value = 2600.919
$ws1.Cells.Item(578,9).Value = 2139.218
$ws1.Cells.Item(578,10).Value = 2816.266
$ws1.Cells.Item(578,12).Value = 8305.8

# Row 579: 2=707.359, 3=1702.27, 4=447.298, 5=1432.259, 6=1912.497, 7=2654.172, 8=2050.083, 9=2055.859, 10=2727.323, 12=7838.347
$ws1.Cells.Item(579,2).Value = 707.359
$ws1.Cells.Item(579,3).Value = 1702.27
$ws1.Cells.Item(579,4).Value = 447.298
$ws1.Cells.Item(579,5).Value = 1432.259
$ws1.Cells.Item(579,6).Value = 1912.497
$ws1.Cells.Item(579,7).Value = 2654.172
$ws1.Cells.Item(579,8).Value = 2050.083
$ws1.Cells.Item(579,9).Value = 2055.859
$ws1.Cells.Item(579,10).Value = 2727.323
$ws1.Cells.Item(579,12).Value = 7838.347

# Row 580: 2=541.988, 3=1467.066, 4=330.048, 5=1190.068, 6=1649.964, 7=2309.591, 8=1568.196, 9=1572.477, 10=2449.006, 12=6531.408
$ws1.Cells.Item(580,2).Value = 541.988
$ws1.Cells.Item(580,3).Value = 1467.066
$ws1.Cells.Item(580,4).Value = 330.048
$ws1.Cells.Item(580,5).Value = 1190.068
$ws1.Cells.Item(580,6).Value = 1649.964
$ws1.Cells.Item(580,7).Value = 2309.591
$ws1.Cells.Item(580,8).Value = 1568.196
$ws1.Cells.Item(580,9).Value = 1572.477
$ws1.Cells.Item(580,10).Value = 2449.006
$ws1.Cells.Item(580,12).Value = 6531.408


# --- Append new row 581 (May 2020) ---
# Copy formatting (date number format) from the row above, then set values.
$ws1.Cells.Item(580,1).Copy() | Out-Null
$ws1.Cells.Item(581,1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws1.Cells.Item(581,1).Value = 43952
$ws1.Cells.Item(581,2).Value = 390.455
$ws1.Cells.Item(581,3).Value = 1451.572
$ws1.Cells.Item(581,4).Value = 257.132
$ws1.Cells.Item(581,5).Value = 1197.248
$ws1.Cells.Item(581,6).Value = 1706.681
$ws1.Cells.Item(581,7).Value = 2423.809
$ws1.Cells.Item(581,8).Value = 1799.352
$ws1.Cells.Item(581,9).Value = 1803.908
$ws1.Cells.Item(581,10).Value = 2722.916
$ws1.Cells.Item(581,11).Value = -6.688
$ws1.Cells.Item(581,12).Value = 6869.847

# --- Update "Annual Data" sheet: 2019 total (row 83, column L) ---
$ws2.Cells.Item(83,12).Value = 100165.861
